# Apply the commit's changes to the deck:
#  1. Three tables (slides 14, 15, 16) switch from the default table
#     style {CAC69132-24B7-49D6-88D9-2FA4395C7D32} to
#     {7B97409D-6A9A-434C-9F9C-10073AFFD174}.
#  2. The slide master's theme (the "Integral" / Red-Violet palette) is
#     swapped for the stock default "Office" palette - i.e. the deck's
#     Design goes from Integral back to the default Office theme colors.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$targetStyle = "{7B97409D-6A9A-434C-9F9C-10073AFFD174}"
for ($idx = 1; $idx -le $p.Slides.Count; $idx++) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($targetStyle)
        }
    }
}

# --- 2. Swap the active theme's color scheme to the Office defaults --
$master = $p.SlideMaster
$scheme = $master.ColorScheme

$scheme.Colors(1).RGB  = 0x000000   # dk1      - 000000
$scheme.Colors(2).RGB  = 0xFFFFFF   # lt1      - FFFFFF
$scheme.Colors(3).RGB  = 0x6A5444   # dk2      - 44546A
$scheme.Colors(4).RGB  = 0xE6E6E7   # lt2      - E7E6E6
$scheme.Colors(5).RGB  = 0xD59B5B   # accent1  - 5B9BD5
$scheme.Colors(6).RGB  = 0x317DED   # accent2  - ED7D31
$scheme.Colors(7).RGB  = 0xA5A5A5   # accent3  - A5A5A5
$scheme.Colors(8).RGB  = 0x00C0FF   # accent4  - FFC000
$scheme.Colors(9).RGB  = 0xC47244   # accent5  - 4472C4
$scheme.Colors(10).RGB = 0x47AD70   # accent6  - 70AD47
$scheme.Colors(11).RGB = 0xC16305   # hlink    - 0563C1
$scheme.Colors(12).RGB = 0x724F95   # folHlink - 954F72
